# "adding averages and more checks"
#
# Training Dashboard:
#  - PERIOD TO EXPIRE (col H) drops by 8 for every data row (one more week-plus
#    has elapsed since the last check).
#  - LAST UPDATE (col I) moves from 08-Sep-2025 to 16-Sep-2025 for every row.
#  - Title + header row get a bold white font (header keeps its dark-blue
#    fill, title loses its old 14pt size).
#
# Exam Dashboard:
#  - COMMENTS (col E) rows 3-6 now all read "date is valid".
#  - COMMENTS column is narrowed (it no longer needs to fit long remarks).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Header / title styling -------------------------------------------------
# Title (merged A1) keeps its bold weight but drops to the default 11pt size
# and turns white.
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").Font.Color = 16777215

$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Color = 16777215

# Header rows turn white-on-blue.
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:E2").Font.Color = 16777215

# --- Training Dashboard: PERIOD TO EXPIRE (H) -------------------------------
$ws1.Range("H3").Value = 399
$ws1.Range("H4").Value = 401
$ws1.Range("H5").Value = 495
$ws1.Range("H6").Value = 403
$ws1.Range("H7").Value = 519
$ws1.Range("H8").Value = 522
$ws1.Range("H9").Value = 520
$ws1.Range("H10").Value = 496
$ws1.Range("H11").Value = 520
$ws1.Range("H12").Value = 699
$ws1.Range("H13").Value = 521
$ws1.Range("H14").Value = 399
$ws1.Range("H15").Value = 524
$ws1.Range("H16").Value = 414
$ws1.Range("H17").Value = 394
$ws1.Range("H18").Value = 524
$ws1.Range("H19").Value = 522
$ws1.Range("H20").Value = 593
$ws1.Range("H21").Value = 82
$ws1.Range("H22").Value = 268
$ws1.Range("H23").Value = 268
$ws1.Range("H24").Value = -36
$ws1.Range("H25").Value = -49
$ws1.Range("H26").Value = 175
$ws1.Range("H27").Value = 215
$ws1.Range("H28").Value = 175
$ws1.Range("H29").Value = 217
$ws1.Range("H30").Value = 216
$ws1.Range("H31").Value = 315
$ws1.Range("H32").Value = 316
$ws1.Range("H33").Value = 318
$ws1.Range("H34").Value = 314
$ws1.Range("H35").Value = 335
$ws1.Range("H36").Value = 336

# --- Training Dashboard: LAST UPDATE (I) ------------------------------------
# These cells hold the date as literal text ("16-Sep-2025"), not a real date
# serial, so go through Formula -> Copy -> PasteSpecial(values) instead of a
# plain .Value assignment (which Excel would otherwise auto-parse into a date
# and reformat the cell).
$ws1.Range("I3:I36").Formula = '="16-Sep-2025"'
$ws1.Range("I3:I36").Copy()
$ws1.Range("I3:I36").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Exam Dashboard: COMMENTS column ----------------------------------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Range("E5").Value = "date is valid"
$ws2.Range("E6").Value = "date is valid"

# Narrow the COMMENTS column now that it just holds short text.
$ws2.Columns.Item(5).ColumnWidth = 14.166666666666666
